# Weekly fruit/vegetable price update:
# Insert a new record as row 325 (pushing the existing rows 325:337 down to
# 326:338) on the "Feria Lagunitas de Puerto Montt - Zapallo italiano" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 325; this shifts rows 325:337
# down to 326:338 and carries the row-above's cell formatting (so D325
# inherits the date number format from D324).
$ws.Rows("325:325").Insert()

# Populate the new row 325 with this week's record.
$ws.Range("A325").Value = 4
$ws.Range("B325").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C325").Value = "Los Lagos"
$ws.Range("D325").Value = 44939
$ws.Range("E325").Value = 10
$ws.Range("F325").Value = 100112032
$ws.Range("G325").Value = "Zapallo italiano"
$ws.Range("H325").Value = "Sin especificar"
$ws.Range("I325").Value = "Primera"
$ws.Range("J325").Value = 240
$ws.Range("K325").Value = 14000
$ws.Range("L325").Value = 15000
$ws.Range("M325").Value = 14500
$ws.Range("N325").Value = "`$/caja 50 unidades"
$ws.Range("O325").Value = "Región de O'Higgins"
$ws.Range("P325").Value = 290
$ws.Range("Q325").Value = 50
$ws.Range("R325").Value = "Hortaliza"
